$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 1; everything currently in row 1 onward shifts down by one.
$ws.Rows.Item(1).Insert()

# Copy the formatting (the bold/bordered/centered header style) that is still sitting on
# what is now row 2 (the old row 1) up onto the brand-new row 1.
$ws.Range("A2:I2").Copy()
$ws.Range("A1:I1").PasteSpecial(-4122)  # -4122 = xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new header row with numeric values 0-8.
for ($i = 0; $i -le 8; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $i
}

# Row 2 (the shifted-down former row 1) should no longer carry the special header
# formatting - it reverts to the default style.
$ws.Range("A2:I2").ClearFormats()

# The former row 1 values for H1/I1 ("thread_size" / "material_surface") are dropped
# when the row shifts down to row 2.
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
